$wb = $excel.ActiveWorkbook

# The "Summary" sheet currently has an extra leading column (A) that holds
# row labels ("Principal", "Interest", "Fees", "Penalties", "Overpaid").
# Remove that whole column so the remaining data shifts left.
$summary = $wb.Worksheets.Item("Summary")
$summary.Columns.Item(1).Delete()

# Select cell B9 on the Summary sheet and make the Summary sheet the active
# (selected) tab of the workbook.
$summary.Activate()
$summary.Range("B9").Select()
